$wb = $excel.ActiveWorkbook

# Update the "想去人数" (number of people interested) figures on both the
# "展览" and "全部类型" sheets for rows 2-5 (F column).
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 1482
    $ws.Range("F3").Value = 3119
    $ws.Range("F4").Value = 46
    $ws.Range("F5").Value = 801
}
